# The deck ships two themes:
#   ppt/theme/theme1.xml -> bound to the (only) Slide Master, currently the
#                            "Integral" / "Red Violet" palette
#   ppt/theme/theme2.xml -> bound to the Notes Master, currently the default
#                            "Office Theme" palette
#
# The authored change swaps the two themes' colour schemes: the slide
# master ends up on the stock "Office" palette and the notes master ends
# up on the "Red Violet" palette that the slide master used to have
# (font scheme / format scheme are identical between the two themes, so
# only the 12 theme colours actually change).
#
# PowerPoint exposes the live theme colours through
#   Master.Theme.ThemeColorScheme.Colors(1..12).RGB
# in a fixed order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation

function Set-ThemeColors {
    param($themeColorScheme, [int[]]$rgbValues)
    for ($i = 1; $i -le $rgbValues.Count; $i++) {
        $themeColorScheme.Colors($i).RGB = $rgbValues[$i - 1]
    }
}

# Target palette for the slide master's theme (theme1.xml): the stock
# "Office" scheme -- dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink.
$officeColors = @(
    0,          # dk1      #000000
    16777215,   # lt1      #FFFFFF
    6968388,    # dk2      #44546A
    15132391,   # lt2      #E7E6E6
    13998939,   # accent1  #5B9BD5
    3243501,    # accent2  #ED7D31
    10855845,   # accent3  #A5A5A5
    49407,      # accent4  #FFC000
    12874308,   # accent5  #4472C4
    4697456,    # accent6  #70AD47
    12673797,   # hlink    #0563C1
    7491477     # folHlink #954F72
)

Set-ThemeColors $p.SlideMaster.Theme.ThemeColorScheme $officeColors
